$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value into a cell without Excel auto-converting
# number-like strings (e.g. "224.24") into a numeric cell. We build the text via
# a formula that evaluates to a string, then Copy + PasteSpecial (paste values) so
# the destination keeps plain-text typing with no extra number formatting/style.
function Set-TextValue($cell, $text) {
    $r = $ws.Range($cell)
    $r.Formula = '="' + $text.Replace('"','""') + '"'
    $r.Copy()
    $r.PasteSpecial(-4163)
}

$ws.Range("D2").Value = "33.846.93"
$ws.Range("E2").Value = "  +9.40%  "

$ws.Range("D3").Value = "1.783.17"
$ws.Range("E3").Value = "  +6.21%  "

$ws.Range("E4").Value = "  +0.02%  "

Set-TextValue 'D5' '224.24'
$ws.Range("E5").Value = "  +2.12%  "

Set-TextValue 'D6' '0.555'
$ws.Range("E6").Value = "  +4.07%  "

$ws.Range("E7").Value = "  +0.08%  "

Set-TextValue 'D8' '30.80'
$ws.Range("E8").Value = "  +5.84%  "

Set-TextValue 'D9' '46.07'
$ws.Range("E9").Value = "  +4.44%  "

Set-TextValue 'D10' '0.278'
$ws.Range("E10").Value = "  +4.92%  "

Set-TextValue 'D11' '0.0657'
$ws.Range("E11").Value = "  +2.44%  "

$ws.Range("E12").Value = "  +1.64%  "

$ws.Range("D13").Value = "2.033.37"
$ws.Range("E13").Value = "  +5.88%  "

$ws.Range("D14").Value = "1.772.81"
$ws.Range("E14").Value = "  +5.95%  "

$ws.Range("E15").Value = "  +3.57%  "

$ws.Range("D16").Value = "33.717.50"
$ws.Range("E16").Value = "  +9.11%  "

Set-TextValue 'D17' '9.94'
$ws.Range("E17").Value = "  -1.13%  "

$ws.Range("E18").Value = "  +1.79%  "

Set-TextValue 'D19' '68.54'
$ws.Range("E19").Value = "  +3.76%  "

Set-TextValue 'D20' '250.82'
$ws.Range("E20").Value = "  +1.75%  "

$ws.Range("D21").Value = "0.0₃0736"
$ws.Range("E21").Value = "  +2.28%  "

Set-TextValue 'D23' '10.23'
$ws.Range("E23").Value = "  +2.64%  "

$ws.Range("E24").Value = "  -0.81%  "

Set-TextValue 'D25' '2.14'
$ws.Range("E25").Value = "  -0.88%  "

Set-TextValue 'D26' '157.66'
$ws.Range("E26").Value = "  -0.97%  "

Set-TextValue 'D27' '16.37'
$ws.Range("E27").Value = "  +3.34%  "

$ws.Range("E28").Value = "  +1.86%  "

$ws.Range("E29").Value = "  +3.11%  "

$ws.Range("E30").Value = "  +0.05%  "

$ws.Range("E31").Value = "  +9.11%  "

Set-TextValue 'D32' '0.0508'
$ws.Range("E32").Value = "  +2.96%  "

$ws.Range("E33").Value = "  +3.90%  "

Set-TextValue 'D34' '3.50'
$ws.Range("E34").Value = "  +5.43%  "

$ws.Range("D35").Value = "1.477.66"
$ws.Range("E35").Value = "  -2.78%  "

Set-TextValue 'D36' '1.74'
$ws.Range("E36").Value = "  -0.30%  "

$ws.Range("E37").Value = "  +2.80%  "

$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue 'D38' '0.619'
$ws.Range("E38").Value = "  +2.41%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue 'D39' '0.0185'
$ws.Range("E39").Value = "  +2.85%  "

Set-TextValue 'D40' '82.63'
$ws.Range("E40").Value = "  -2.18%  "

$ws.Range("B41").Value = "HuobiToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue 'D41' '2.36'
$ws.Range("E41").Value = "  +2.84%  "

$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue 'D42' '2.72'
$ws.Range("E42").Value = "  +2.73%  "

$ws.Range("E43").Value = "  +4.95%  "

$ws.Range("E44").Value = "  +0.45%  "

Set-TextValue 'D45' '0.0506'
$ws.Range("E45").Value = "  +0.86%  "

$ws.Range("E46").Value = "  +3.75%  "

$ws.Range("D47").Value = "1.924.76"
$ws.Range("E47").Value = "  +6.01%  "

$ws.Range("E48").Value = "  +0.21%  "

Set-TextValue 'D49' '5.67'
$ws.Range("E49").Value = "  +1.83%  "

Set-TextValue 'D50' '11.79'
$ws.Range("E50").Value = "  +11.15%  "

Set-TextValue 'D51' '50.70'
$ws.Range("E51").Value = "  -0.61%  "
